$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Nombre de la Empresa" values in column A for the remaining rows
$ws.Range("A2").Value = "CAMPERO"
$ws.Range("A3").Value = "CAMPERO2"
$ws.Range("A4").Value = "CAMPERO23"

# "Años de trayectoria" for row 2 changes from 20 to 202
$ws.Range("C2").Value = 202

# "Email" column value changes from 45454 to 45454sdsdsd for all remaining rows
$ws.Range("E2").Value = "45454sdsdsd"
$ws.Range("E3").Value = "45454sdsdsd"
$ws.Range("E4").Value = "45454sdsdsd"

# Remove the now-obsolete rows 5 and 6
$ws.Rows("5:6").Delete()

# Set column K (11th column) width to 30
# Note: the stored OOXML "width" attribute is ColumnWidth + 5px padding
# (~0.8333 chars at the default font), so back that padding out here so
# the saved width attribute comes out to exactly 30.
$ws.Columns(11).ColumnWidth = 29.1666666666667
